# Add "hydrogen combined cycle" as a new power plant type row on the
# FSPbPPT sheet, and rename the existing "hydrogen" row to
# "hydrogen combustion turbine".

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("FSPbPPT")

# Rename existing row 24 label from "hydrogen" to "hydrogen combustion turbine"
# and give it the black-font / vertically-centered style used for the new row.
$ws.Range("A24").Value = "hydrogen combustion turbine"
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108

# Add the new "hydrogen combined cycle" row, copying A24's style for the label
# cell and B24's style for the value cell.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("B25").Value = 0
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Leave the cell selection on the new area, like the author's last
# position, but keep the "About" sheet as the active tab of the workbook.
$ws.Range("B31").Select()
$wsAbout.Activate()
